# Updated capital structure database
# Refreshes the Venezuela "Bank (Money Center)" capital-structure rows:
#  - Row 2 (industry aggregate "4") gets refreshed metric values.
#  - Rows 3-5 company metrics are refreshed, and the company names are
#    rotated one row down (Banco Provincial -> row4, Banco Nacional de
#    Credito -> row5, Mercantil Servicios Financieros -> row3).
#  - Row 6 (Banco del Caribe) gets refreshed metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 - "4" (industry rollup row)
# ---------------------------------------------------------------------
$ws.Range("D2").Value = 1.9445
$ws.Range("E2").Value = 1.575
$ws.Range("K2").Value = 58.082
$ws.Range("L2").Value = 0.4313233328382594
$ws.Range("M2").Value = 0.006
$ws.Range("N2").Value = 0.0000008178286648947045
$ws.Range("O2").Value = 0.0001033022278847147
$ws.Range("P2").Value = 0.006
$ws.Range("Q2").Value = 0.0000008178286648947045
$ws.Range("R2").Value = 0.0001033022278847147
$ws.Range("U2").Value = 347.9
$ws.Range("V2").Value = 0.04742043208614462
$ws.Range("W2").Value = 0.06803840452838764
$ws.Range("X2").Value = 0.1238511038859222
$ws.Range("Y2").Value = -0.05581269935753456
$ws.Range("Z2").Value = -0.982131135584567
$ws.Range("AB2").Value = 0.1238517623429393
$ws.Range("AC2").Value = -0.1238517623429393
$ws.Range("AD2").Value = 1.346
$ws.Range("AF2").Value = 1.346
$ws.Range("AG2").Value = -346.554
$ws.Range("AH2").Value = 0.0001834325768079625
$ws.Range("AI2").Value = 0.006983646023576292
$ws.Range("AJ2").Value = -0.04957892378567731
$ws.Range("AK2").Value = 2.233469103658065

# ---------------------------------------------------------------------
# Row 3 - now "Mercantil Servicios Financieros, C.A. (CCSE:MVZ.B)"
# ---------------------------------------------------------------------
$ws.Range("B3").Value = "Mercantil Servicios Financieros, C.A. (CCSE:MVZ.B)"
$ws.Range("D3").Value = 2.575
$ws.Range("E3").Value = 2.981
$ws.Range("K3").Value = 54.1
$ws.Range("L3").Value = 0.4967860422405877
$ws.Range("M3").Value = 0.001
$ws.Range("N3").Value = 0.000002447381302006853
$ws.Range("O3").Value = 0.00001848428835489834
$ws.Range("P3").Value = 0.001
$ws.Range("Q3").Value = 0.000002447381302006853
$ws.Range("R3").Value = 0.00001848428835489834
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 128.4
$ws.Range("V3").Value = 0.3142437591776799
$ws.Range("W3").Value = 0.6320093457943926
$ws.Range("X3").Value = 0.1238467486970918
$ws.Range("Y3").Value = 0.5081625970973008
$ws.Range("Z3").Value = -1.930851063829787
$ws.Range("AB3").Value = 0.1238467486970918
$ws.Range("AC3").Value = -0.1238467486970918
$ws.Range("AG3").Value = -128.4
$ws.Range("AJ3").Value = -0.4582441113490364
$ws.Range("AK3").Value = 2.562874251497006

# ---------------------------------------------------------------------
# Row 4 - now "Banco Provincial, S.A. Banco Universal (CCSE:BPV)"
# ---------------------------------------------------------------------
$ws.Range("B4").Value = "Banco Provincial, S.A. Banco Universal (CCSE:BPV)"
$ws.Range("D4").Value = 1.43
$ws.Range("E4").Value = 0.9359999999999999
$ws.Range("K4").Value = 1.7
$ws.Range("L4").Value = 0.1240875912408759
$ws.Range("U4").Value = 138.9
$ws.Range("V4").Value = 0.8170588235294118
$ws.Range("W4").Value = 0.02407932011331445
$ws.Range("X4").Value = 0.1238467486970918
$ws.Range("Y4").Value = -0.09976742858377738
$ws.Range("Z4").Value = -0.2264462809917355
$ws.Range("AB4").Value = 0.1238467486970918
$ws.Range("AC4").Value = -0.1238467486970918
$ws.Range("AG4").Value = -138.9
$ws.Range("AJ4").Value = -4.466237942122188
$ws.Range("AK4").Value = 2.225961538461538

# ---------------------------------------------------------------------
# Row 5 - now "Banco Nacional de Crédito, C.A., Banco Universal (CCSE:BNC)"
# ---------------------------------------------------------------------
$ws.Range("B5").Value = "Banco Nacional de Crédito, C.A., Banco Universal (CCSE:BNC)"
$ws.Range("D5").Value = 2.021
$ws.Range("E5").Value = 2.08
$ws.Range("K5").Value = 1.78
$ws.Range("L5").Value = 0.2772585669781932
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("T5").ClearContents()
$ws.Range("U5").Value = 67.59999999999999
$ws.Range("V5").Value = 0.01015457181054813
$ws.Range("W5").Value = 0.06425992779783393
$ws.Range("X5").Value = 0.1238554590747526
$ws.Range("Y5").Value = -0.05959553127691865
$ws.Range("Z5").Value = -0.3258883248730964
$ws.Range("AA5").Value = -0
$ws.Range("AB5").Value = 0.1238567759887867
$ws.Range("AC5").Value = -0.1238567759887867
$ws.Range("AD5").Value = 0.767
$ws.Range("AF5").Value = 0.767
$ws.Range("AG5").Value = -66.833
$ws.Range("AH5").Value = 0.0001152020609603646
$ws.Range("AI5").Value = 0.02559482096973337
$ws.Range("AJ5").Value = -0.01014116726985416
$ws.Range("AK5").Value = 1.775914755666569

# ---------------------------------------------------------------------
# Row 6 - "Banco del Caribe, C.A., Banco Universal (CCSE:ABC.A)" (unchanged name)
# ---------------------------------------------------------------------
$ws.Range("D6").Value = 1.868
$ws.Range("E6").Value = 1.07
$ws.Range("K6").Value = 0.502
$ws.Range("L6").Value = 0.08900709219858156
$ws.Range("M6").Value = 0.005
$ws.Range("N6").Value = 0.0000496031746031746
$ws.Range("O6").Value = 0.009960159362549801
$ws.Range("P6").Value = 0.005
$ws.Range("Q6").Value = 0.0000496031746031746
$ws.Range("R6").Value = 0.009960159362549801
$ws.Range("U6").Value = 13
$ws.Range("V6").Value = 0.128968253968254
$ws.Range("W6").Value = 0.07181688125894134
$ws.Range("X6").Value = 0.1242810036033131
$ws.Range("Y6").Value = -0.05246412234437177
$ws.Range("Z6").Value = -11.05882352941177
$ws.Range("AB6").Value = 0.1243438603696422
$ws.Range("AC6").Value = -0.1243438603696422
$ws.Range("AD6").Value = 0.579
$ws.Range("AF6").Value = 0.579
$ws.Range("AG6").Value = -12.421
$ws.Range("AH6").Value = 0.005711241973189714
$ws.Range("AI6").Value = 0.07265654410842012
$ws.Range("AJ6").Value = -0.1405424365516695
$ws.Range("AK6").Value = 2.468892864241702
